# Update cfb_weather.xlsx with Timestamp 2024-11-02T10:01:36.050813
#
# 1. Refresh the run timestamp stamped down column AK ("Timestamp") of the
#    "FBS" sheet (rows 2-45) from 2024-11-02T05:15:55.011945 to
#    2024-11-02T10:01:36.050813.
# 2. Correct a handful of forecast wind-direction ("wind_dir_fg") readings
#    that were re-pulled since the timestamp moved:
#      FBS!Q16  SW  -> SSW
#      FBS!Q18  N   -> E
#      FBS!Q32  SW  -> SSW
#      FBS!Q34  NNW -> NW
#      Other!S4  NNE -> N
#      Other!S13 SE  -> SSE

$wb = $excel.ActiveWorkbook

$fbs = $wb.Worksheets.Item("FBS")
$other = $wb.Worksheets.Item("Other")

$newTimestamp = "2024-11-02T10:01:36.050813"

for ($row = 2; $row -le 45; $row++) {
    $fbs.Range("AK$row").Value = $newTimestamp
}

$fbs.Range("Q16").Value = "SSW"
$fbs.Range("Q18").Value = "E"
$fbs.Range("Q32").Value = "SSW"
$fbs.Range("Q34").Value = "NW"

$other.Range("S4").Value = "N"
$other.Range("S13").Value = "SSE"
